# Refresh "想去人数" (interest count, column F) figures on each sheet to the
# values captured in the latest scrape (gh-pages data regeneration).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 625
$ws.Range("F3").Value = 626
$ws.Range("F4").Value = 910
$ws.Range("F5").Value = 674
$ws.Range("F6").Value = 804
$ws.Range("F7").Value = 373
$ws.Range("F8").Value = 577
$ws.Range("F9").Value = 117
$ws.Range("F10").Value = 1168
$ws.Range("F11").Value = 598
$ws.Range("F12").Value = 355
$ws.Range("F13").Value = 479
$ws.Range("F14").Value = 154
$ws.Range("F15").Value = 80
$ws.Range("F16").Value = 315
$ws.Range("F18").Value = 75
$ws.Range("F19").Value = 532
$ws.Range("F20").Value = 42
$ws.Range("F21").Value = 544
$ws.Range("F23").Value = 561
$ws.Range("F24").Value = 2

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 58
$ws.Range("F6").Value = 17
$ws.Range("F8").Value = 173
$ws.Range("F13").Value = 45

# Sheet "全部类型" (All types, a combined/aggregated view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 625
$ws.Range("F5").Value = 58
$ws.Range("F7").Value = 626
$ws.Range("F8").Value = 910
$ws.Range("F9").Value = 674
$ws.Range("F10").Value = 804
$ws.Range("F11").Value = 373
$ws.Range("F12").Value = 577
$ws.Range("F13").Value = 117
$ws.Range("F14").Value = 1168
$ws.Range("F15").Value = 598
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 355
$ws.Range("F19").Value = 479
$ws.Range("F21").Value = 154
$ws.Range("F22").Value = 80
$ws.Range("F23").Value = 173
$ws.Range("F24").Value = 315
$ws.Range("F26").Value = 75
$ws.Range("F29").Value = 532
$ws.Range("F32").Value = 45
$ws.Range("F33").Value = 42
$ws.Range("F34").Value = 544
$ws.Range("F36").Value = 561
$ws.Range("F37").Value = 2

